# Generate Report for Handoff
#
# Updates the localization-status report to reflect a fresh handoff run:
#   - Status text moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   - The associated timestamps advance a few seconds (new handoff generated)
#   - The (previously over-wide) status-column widths shrink to fit the new,
#     shorter "Ready for handoff" text
#
# The workbook has three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# E2 / F2 = status for zh-cn / de-de, G2 = latest HO xliff generate date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-04 07:05:46"

# Status columns (E, F) are no longer the widest text on the sheet now
# that they read "Ready for handoff" instead of "Handed back: in sync
# with en-US" -- narrow them to match.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet -------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-04 07:05:42"
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de sheet -------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-04 07:05:46"
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
